$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Change 2 (done first): remove the "_GoBack" bookmark that currently sits
# between the "()" run and the following single-space run, and fold that
# space into the "()" run so the text reads "() " as one run, leaving the
# "%}" run that follows untouched and separate.
# ---------------------------------------------------------------------------
$oldBookmark = $d.Bookmarks("_GoBack")
$oldStart = $oldBookmark.Start

$spaceRange = $d.Range($oldStart, $oldStart + 1)
$spaceRange.Delete()

$parenRange = $d.Range($oldStart - 2, $oldStart)
$parenRange.InsertAfter(" ")

$d.Bookmarks("_GoBack").Delete()

# ---------------------------------------------------------------------------
# Change 1: move the "_GoBack" bookmark to the very start of the document
# (right after the first paragraph's pPr, before its first run).
#
# A zero-length Range(0,0) at the literal start of the body places the
# bookmarkStart/bookmarkEnd pair in two different spots, so instead we
# insert a throw-away character at position 0, anchor the bookmark right
# after it (a non-zero offset, which behaves correctly), and then delete
# the throw-away character again. Word keeps the (now zero-length)
# bookmark anchored at that same spot.
# ---------------------------------------------------------------------------
$startRange = $d.Range(0, 0)
$startRange.InsertBefore("X")

$afterMarker = $d.Range(1, 1)
$d.Bookmarks.Add("_GoBack", $afterMarker)

$markerRange = $d.Range(0, 1)
$markerRange.Delete()
